# Refresh the "cryptos" price/volume snapshot (scheduled GitHub Actions scrape).
# For each changed row: columns B/C (coin name/link) are only touched for the
# Cosmos/Kaspa row swap (rows 42-43); columns D (price) and E (1h volume %) are
# refreshed with the latest scraped figures for every listed coin.
#
# Column D holds prices as plain TEXT (the source sheet never stores them as
# numbers - e.g. "64.243.88" is not even a valid number). Excel auto-converts a
# plain numeric-looking string typed into .Value to a Number, so any new D value
# that parses as a float is written with a leading apostrophe to force text,
# exactly like typing  `'0.476`  into Excel by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "64.223.18"
$ws.Range("E2").Value = "  -3.57%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.155.05"
$ws.Range("E3").Value = "  -2.34%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5: BNB
$ws.Range("D5").Value = "'607.30"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6: Solana
$ws.Range("D6").Value = "'146.17"
$ws.Range("E6").Value = "  -6.33%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8: LidoStakedEther
$ws.Range("D8").Value = "3.151.44"
$ws.Range("E8").Value = "  -2.40%  "

# Row 9: XRP
$ws.Range("E9").Value = "  -3.46%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -6.41%  "

# Row 11: Toncoin
$ws.Range("E11").Value = "  -4.70%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.476"
$ws.Range("E12").Value = "  -5.00%  "

# Row 13: ShibaInu
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -6.02%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'36.08"
$ws.Range("E14").Value = "  -7.14%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.673.50"
$ws.Range("E15").Value = "  -2.14%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "64.226.32"
$ws.Range("E16").Value = "  -3.55%  "

# Row 17: TRON
$ws.Range("E17").Value = "  +1.18%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.159.38"
$ws.Range("E18").Value = "  -3.30%  "

# Row 19: Polkadot
$ws.Range("E19").Value = "  -4.39%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'478.72"
$ws.Range("E20").Value = "  -5.62%  "

# Row 21: Chainlink
$ws.Range("D21").Value = "'14.61"
$ws.Range("E21").Value = "  -4.19%  "

# Row 22: Polygon
$ws.Range("E22").Value = "  -3.69%  "

# Row 23: Uniswap
$ws.Range("D23").Value = "'7.73"
$ws.Range("E23").Value = "  -3.57%  "

# Row 24: InternetComputer(DFINITY)
$ws.Range("E24").Value = "  -5.50%  "

# Row 25: Litecoin
$ws.Range("D25").Value = "'82.96"
$ws.Range("E25").Value = "  -3.69%  "

# Row 26: Dai
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "

# Row 27: PancakeSwap
$ws.Range("D27").Value = "'2.89"
$ws.Range("E27").Value = "  -3.47%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "'8.44"
$ws.Range("E28").Value = "  -6.80%  "

# Row 29: ImmutableX
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -6.69%  "

# Row 30: Hedera
$ws.Range("D30").Value = "'0.118"
$ws.Range("E30").Value = "  -28.00%  "

# Row 31: NEARProtocol
$ws.Range("D31").Value = "'6.84"
$ws.Range("E31").Value = "  -2.00%  "

# Row 32: Stacks
$ws.Range("D32").Value = "'2.75"
$ws.Range("E32").Value = "  -5.14%  "

# Row 33: FirstDigitalUSD
$ws.Range("E33").Value = "  +0.07%  "

# Row 34: EthereumClassic
$ws.Range("D34").Value = "'26.21"
$ws.Range("E34").Value = "  -7.00%  "

# Row 35: Mantle
$ws.Range("E35").Value = "  -5.40%  "

# Row 36: Filecoin
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  -5.55%  "

# Row 37: OKB
$ws.Range("D37").Value = "'54.08"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38: PEPE
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("E38").Value = "  -8.59%  "

# Row 39: Bittensor
$ws.Range("D39").Value = "'451.46"
$ws.Range("E39").Value = "  -8.36%  "

# Row 40: dogwifhat
$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = "  -6.18%  "

# Row 41: VeChain
$ws.Range("D41").Value = "'0.0397"
$ws.Range("E41").Value = "  -5.76%  "

# Row 42: Cosmos
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  -7.22%  "

# Row 43: Kaspa
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.46"
$ws.Range("E43").Value = "  -2.87%  "

# Row 44: Maker
$ws.Range("D44").Value = "2.848.48"
$ws.Range("E44").Value = "  -3.24%  "

# Row 45: TheGraph
$ws.Range("D45").Value = "'0.268"
$ws.Range("E45").Value = "  -8.33%  "

# Row 46: Fetch.AI
$ws.Range("E46").Value = "  -8.20%  "

# Row 47: InjectiveProtocol
$ws.Range("E47").Value = "  -6.43%  "

# Row 48: USDe
$ws.Range("D48").Value = "'0.998"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49: ThetaToken
$ws.Range("E49").Value = "  -4.02%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  -3.31%  "

# Row 51: Monero
$ws.Range("D51").Value = "'118.57"
$ws.Range("E51").Value = "  -2.02%  "
